# modified test cases on overdue fix
# Recompute of the loan repayment schedule after an "overdue" fix: a few
# rows shift by a cent or two and an extra (7th) instalment row appears
# that finally drives the loan balance to zero.

$wb  = $excel.ActiveWorkbook
$wsSummary    = $wb.Worksheets.Item("Summary")
$wsSchedule   = $wb.Worksheets.Item("Repayment schedule")
$wsTxn        = $wb.Worksheets.Item("Transactions")

# ---------------------------------------------------------------------
# 1) Repayment schedule sheet - re-normalise the number formats that had
#    drifted onto a stray "0.00" style back to the sheet's standard
#    general/vertical-center/wrap style (and the couple of special
#    formats used elsewhere on the sheet), reusing existing styles so we
#    don't create new ones.
# ---------------------------------------------------------------------

# Date style for the Date column
$wsSchedule.Range("C3").Copy() | Out-Null
$wsSchedule.Range("C2:C9").PasteSpecial(-4122) | Out-Null

# Disbursement amount cell keeps the thousands style used on Summary!A2
$wsSummary.Range("A2").Copy() | Out-Null
$wsSchedule.Range("G2").PasteSpecial(-4122) | Out-Null

# Balance-of-loan column (rows 3-6) keeps the 2 decimal style used on Summary!F2
$wsSummary.Range("F2").Copy() | Out-Null
$wsSchedule.Range("G3:G6").PasteSpecial(-4122) | Out-Null

# Everything else in the data block reverts to the plain general style
$wsSchedule.Range("A2").Copy() | Out-Null
$wsSchedule.Range("A2:B9").PasteSpecial(-4122) | Out-Null
$wsSchedule.Range("D2:F9").PasteSpecial(-4122) | Out-Null
$wsSchedule.Range("H2:P9").PasteSpecial(-4122) | Out-Null
$wsSchedule.Range("G7:G9").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 2) Summary sheet - Over Due column picks up the corrected figure and
#    also reverts to the plain general style.
# ---------------------------------------------------------------------
$wsSummary.Range("B3").Copy() | Out-Null
$wsSummary.Range("A3").PasteSpecial(-4122) | Out-Null
$wsSummary.Range("E3").PasteSpecial(-4122) | Out-Null
$wsSummary.Range("F3").PasteSpecial(-4122) | Out-Null

$wsSummary.Range("A3").Value = 116.11
$wsSummary.Range("E3").Value = 116.11
$wsSummary.Range("F3").Value = 40

# ---------------------------------------------------------------------
# 3) Repayment schedule values - clear the stray P2 total, correct the
#    rounding on a handful of existing rows and fill in the brand new
#    7th instalment row that finally zeroes the loan out.
# ---------------------------------------------------------------------
$wsSchedule.Range("P2").Value = ""

$wsSchedule.Range("G5").Value = 2520.2199999999998

$wsSchedule.Range("F6").Value = 830.54
$wsSchedule.Range("G6").Value = 1689.68

$wsSchedule.Range("G7").Value = 852.78

$wsSchedule.Range("F8").Value = 844.24
$wsSchedule.Range("G8").Value = 8.5399999999999991
$wsSchedule.Range("P8").Value = 851.25

$wsSchedule.Range("A9").Value = 7
$wsSchedule.Range("B9").Value = 31
$wsSchedule.Range("C9").Value = 42217
$wsSchedule.Range("F9").Value = 8.5399999999999991
$wsSchedule.Range("G9").Value = 0
$wsSchedule.Range("H9").Value = 0.07
$wsSchedule.Range("I9").Value = 0
$wsSchedule.Range("J9").Value = 0
$wsSchedule.Range("K9").Value = 8.61
$wsSchedule.Range("L9").Value = 0
$wsSchedule.Range("M9").Value = 0
$wsSchedule.Range("N9").Value = 0
$wsSchedule.Range("O9").Value = 0
$wsSchedule.Range("P9").Value = 8.61

# ---------------------------------------------------------------------
# 4) Re-create the cursor position each sheet was left on, finishing
#    back on the Transactions tab (which stays the active tab).
# ---------------------------------------------------------------------
$wsSummary.Activate() | Out-Null
$wsSummary.Range("A7").Select() | Out-Null

$wsSchedule.Activate() | Out-Null
$wsSchedule.Range("A11").Select() | Out-Null

$wsTxn.Activate() | Out-Null
$wsTxn.Range("A5").Select() | Out-Null
